# Schedule update: classes that span multiple 30-minute slots are now
# represented by a single merged cell (anchored at the slot where the class
# starts) showing the course name plus its start/end time, instead of being
# repeated in every slot row it occupies. Row heights are also reduced from
# 50 to 20 now that most slot rows are empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rewrite cell contents for the slots that changed ---------------

$newValues = @{
    "B2"  = $null
    "C2"  = $null
    "D2"  = $null
    "E2"  = "EOSC_V 111-L1E - Laboratory Exploration of Planet Earth`n8:00 a.m. - 11:00 a.m."
    "F2"  = $null

    "B3"  = $null
    "C3"  = $null
    "E3"  = $null

    "B4"  = "CPSC_V 221-L1J - Basic Algorithms and Data Structures`n9:00 a.m. - 11:00 a.m."
    "E4"  = $null

    "E5"  = $null

    "C8"  = "MATH_V 200-102 - Calculus III`n11:00 a.m. - 12:30 p.m."
    "E8"  = "MATH_V 200-102 - Calculus III`n11:00 a.m. - 12:30 p.m."

    "E12" = "CPSC_V 330-T1F - Applied Machine Learning`n1:00 p.m. - 2:00 p.m."

    "B14" = "CPSC_V 221-101 - Basic Algorithms and Data Structures`n2:00 p.m. - 3:00 p.m."
    "D14" = "CPSC_V 221-101 - Basic Algorithms and Data Structures`n2:00 p.m. - 3:00 p.m."
    "F14" = "CPSC_V 221-101 - Basic Algorithms and Data Structures`n2:00 p.m. - 3:00 p.m."

    "C20" = "CPSC_V 330-103 - Applied Machine Learning`n5:00 p.m. - 6:30 p.m."
    "E20" = "CPSC_V 330-103 - Applied Machine Learning`n5:00 p.m. - 6:30 p.m."
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}

# --- 2. Merge the cells that now represent a multi-slot class ----------

$mergeRanges = @(
    "E2:E7",
    "B4:B7",
    "C8:C10",
    "E8:E10",
    "E12:E13",
    "B14:B15",
    "D14:D15",
    "F14:F15",
    "C20:C22",
    "E20:E22"
)

foreach ($rng in $mergeRanges) {
    $ws.Range($rng).Merge()
}

# --- 3. Shrink the slot rows now that most of them are empty -----------

for ($r = 1; $r -le 29; $r++) {
    $ws.Rows.Item($r).RowHeight = 20
}

# --- 4. Tidy up the trailing 10:00 p.m. row -----------------------------
# Drop the (now unused) formatting on A30 and fully clear B30:F30, which
# were never used by any class. (Use Value2 to read back the label -- the
# Value getter does not resolve cleanly when read as a plain expression.)

$lastLabel = $ws.Range("A30").Value2
$ws.Range("A30:F30").Clear()
$ws.Range("A30").Value = $lastLabel
